# Method of Hierarchy Analysis
# - Wrap the long model/header names onto multiple lines (space -> newline)
#   in the header row (B1:H1) and in the row labels (A2:A7).
# - Recompute the priority-vector values in column H (H2, H3, H5, H7).
# - Shrink the column widths now that the headers wrap onto several lines.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10

# ---- Header row (row 1): car/model names now wrap onto separate lines ----
$ws.Range("B1").Value = "Kia" + $nl + "Rio"
$ws.Range("C1").Value = "Volkswagen" + $nl + "Golf"
$ws.Range("D1").Value = "Toyota" + $nl + "Corolla"
$ws.Range("E1").Value = "Skoda" + $nl + "Octavia"
$ws.Range("F1").Value = "BMW" + $nl + "3" + $nl + "Series"
$ws.Range("G1").Value = "Hyundai" + $nl + "Solaris"
$ws.Range("H1").Value = "Вектор" + $nl + "приоритетов"

# ---- Row labels (column A), same wrapping applied ----
$ws.Range("A2").Value = "Kia" + $nl + "Rio"
$ws.Range("A3").Value = "Volkswagen" + $nl + "Golf"
$ws.Range("A4").Value = "Toyota" + $nl + "Corolla"
$ws.Range("A5").Value = "Skoda" + $nl + "Octavia"
$ws.Range("A6").Value = "BMW" + $nl + "3" + $nl + "Series"
$ws.Range("A7").Value = "Hyundai" + $nl + "Solaris"

# ---- Updated priority-vector values ----
$ws.Range("H2").Value = "0.222"
$ws.Range("H3").Value = "0.092"
$ws.Range("H5").Value = "0.222"
$ws.Range("H7").Value = "0.092"

# ---- Narrower column widths to suit the now-wrapped, shorter header lines ----
$ws.Columns.Item(1).ColumnWidth = 13.5167
$ws.Columns.Item(2).ColumnWidth = 5.1167
$ws.Columns.Item(3).ColumnWidth = 13.5167
$ws.Columns.Item(4).ColumnWidth = 9.9167
$ws.Columns.Item(5).ColumnWidth = 9.9167
$ws.Columns.Item(6).ColumnWidth = 8.75
$ws.Columns.Item(7).ColumnWidth = 9.9167
$ws.Columns.Item(8).ColumnWidth = 14.75
